# Update "想去人数" (interested-count) figures across sheets, as generated
# at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value  = 369
$ws.Range("F4").Value  = 420
$ws.Range("F5").Value  = 1149
$ws.Range("F8").Value  = 984
$ws.Range("F10").Value = 6115
$ws.Range("F12").Value = 1767
$ws.Range("F14").Value = 6020
$ws.Range("F15").Value = 6020
$ws.Range("F21").Value = 839
$ws.Range("F24").Value = 1433
$ws.Range("F26").Value = 265
$ws.Range("F31").Value = 3867

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value  = 314
$ws.Range("F5").Value  = 174
$ws.Range("F8").Value  = 385
$ws.Range("F12").Value = 5

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value  = 9508
$ws.Range("F4").Value  = 636
$ws.Range("F5").Value  = 206

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value  = 9508
$ws.Range("F4").Value  = 636
$ws.Range("F5").Value  = 369
$ws.Range("F6").Value  = 420
$ws.Range("F7").Value  = 1149
$ws.Range("F11").Value = 314
$ws.Range("F12").Value = 984
$ws.Range("F13").Value = 206
$ws.Range("F15").Value = 6115
$ws.Range("F17").Value = 1767
$ws.Range("F22").Value = 5
$ws.Range("F23").Value = 6020
$ws.Range("F24").Value = 6020
$ws.Range("F30").Value = 839
$ws.Range("F33").Value = 1433
$ws.Range("F36").Value = 265
$ws.Range("F46").Value = 3867
